# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计", with the
#   quarter's per-fund holding detail (same shape as the "2021-Q4" sheet).
# - Insert a new leading row into the "总计" (totals) sheet summarising the
#   2022-Q1 totals (8 funds, 6.49 billion yuan held), pushing the older
#   quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the "2022-Q1" sheet right after "2021-Q4" (before "总计").
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4)
$newSheet.Name = "2022-Q1"

# Seed it from "2021-Q4" so it inherits the same header/column styling
# (bold, centered, bordered header row + first-column style), then overwrite
# every value below.
$q4.Range("A1:H9").Copy($newSheet.Range("A1"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund codes (e.g. "010108") and the D:G metrics are stored as *text* in the
# source data (leading zeros / fixed decimal places matter), so force the
# text number format before writing them - otherwise "010108" would be
# auto-coerced to the number 10108.
$newSheet.Range("B2:G9").NumberFormat = "@"

$rows = @(
    @("010108", "景顺长城核心招景混合", "54.59", "89.90", "4.00", "2.1836", 7),
    @("010027", "景顺长城核心中景一年持有期混合", "53.17", "90.70", "4.10", "2.1800", 7),
    @("009190", "景顺长城核心优选一年持有期混合", "18.91", "89.60", "10.37", "1.9610", 1),
    @("008060", "景顺长城价值边际灵活配置混合", "4.93", "80.78", "1.99", "0.0981", 10),
    @("010783", "德邦沪港深龙头混合A", "0.93", "81.58", "5.40", "0.0502", 4),
    @("010784", "德邦沪港深龙头混合C", "0.27", "81.58", "5.40", "0.0146", 4),
    @("519602", "海富通大中华精选混合QDII", "0.11", "89.68", "3.99", "0.0044", 8),
    @("160922", "大成恒生综合中小型股指数(QDII-LOF)A", "0.10", "92.44", "1.56", "0.0016", 4)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).Value = $row[4]
    $newSheet.Cells.Item($r, 7).Value = $row[5]
    $newSheet.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Insert the 2022-Q1 summary row at the top of the "总计" sheet.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows(2).Insert()

# The inserted row picked up the header row's formatting (border/bold) on
# B2:D2 - clear that back to the plain style the data rows use, then give A2
# the same style as the existing data rows (A3 below, copied down with the
# insert) so it matches A3/A4.
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 8
$total.Range("D2").Value = 6.49

Write-Host "2022-Q1 sheet added and 总计 updated"
